$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.205.25'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.36%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.905.98'
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.27'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.46%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.22%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5151'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.74%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4018'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.27%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08470'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.50%  '
$ws.Range("E10").Value = '  +0.26%  '
$ws.Range("E11").Value = '  -0.20%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '23.26'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +13.06%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.443'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.14%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.906.00'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.43%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.357'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.10%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.28%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '94.83'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.88%  '
$ws.Range("E18").Value = '  +0.32%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06685'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.49%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.35'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.65%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.002'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.39%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '30.209.71'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.39%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.27'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.77%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.205'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.54%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.125.65'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.51%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.64'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.15%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '161.43'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.91%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.387'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.85%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '129.50'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.099'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.50%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1058'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.36%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.049'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.04%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.705'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.59%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02496'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.82%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06569'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.71%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2200'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.35%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.196'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.83%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.228'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.18%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '11.90'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.67%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.809'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.87%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6514'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.78%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.231'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.43%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6123'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.46%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.24'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.38%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.717'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.16%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.062'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.80%  '
$ws.Range("E48").Value = '  +1.17%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '125.00'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.56%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.160'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.57%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.22'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.78%  '
